# The document's header/footer each carry a small logo inline picture.
# This edit simply swaps the "name" each logo picture is tagged with
# (the wp:docPr / pic:cNvPr "name" attribute - cosmetic metadata, not the
# embedded media filename or the alt-text/description):
#
#   header -> BTec_Logo-Orange picture : image1.jpg -> image2.jpg
#   footer -> PearsonLogo picture (x2) : image2.png  -> image1.png
#
# NB: setting InlineShape.Name directly on a shape obtained straight from
# a Header/Footer Range silently fails to stick in this runtime, so each
# shape's range is selected first and the rename is applied through
# $word.Selection.InlineShapes - that path persists correctly for both
# header- and footer-hosted pictures.

$d = $word.ActiveDocument

function Rename-LogoInlineShape($range, $oldName, $newName) {
    if ($range.InlineShapes.Count -le 0) {
        return $false
    }
    $shape = $range.InlineShapes.Item(1)
    if ($shape.Name -ne $oldName) {
        return $false
    }
    $shape.Range.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
    return $true
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    for ($hi = 1; $hi -le 3; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            Rename-LogoInlineShape $hdr.Range "image1.jpg" "image2.jpg" | Out-Null
        }

        $ftr = $sec.Footers.Item($hi)
        if ($ftr.Exists) {
            Rename-LogoInlineShape $ftr.Range "image2.png" "image1.png" | Out-Null
        }
    }
}

Write-Output "Renamed logo inline pictures."
